$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 1, shifting existing data down.
$ws.Rows.Item(1).Insert()

# New row 1: flow-velocity style formulas (Cylinder Table 1 Data)
$ws.Range("A1").Formula = "=20/328"
$ws.Range("B1").Formula = "=20/328"
$ws.Range("C1").Formula = "=20/135.75"
$ws.Range("D1").Formula = "=20/135.75"
$ws.Range("E1").Formula = "=20/133.5"
$ws.Range("F1").Formula = "=20/133.5"

$ws.Range("H1").Value = "V_flow"
$ws.Range("I1").Value = 132000

# Match the selection recorded in the saved file
$ws.Range("M10").Select()

Write-Output "Applied Cylinder Table 1 Data edits"
